# Power BI Technical Test
#
# Applies the target edits to the Sharepoint Form answers:
#  1. Rewrite the "Power Apps experience" answer.
#  2. Rewrite the "weakness" answer.
#  3. Split the "Why do you want to work at Data-Driven?" heading so it is
#     immediately followed by a brand-new answer paragraph, and move the
#     (rendering-cache) page-break hint onto that heading.
#  4. Drop the page-break hint that used to sit on the visa-answer
#     paragraph (it now lives on the "Why..." heading instead).

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $startText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.StartsWith($startText)) {
            return $i
        }
    }
    return -1
}

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. "If you have any experience with Power Apps" answer -----------------
$idx = Find-ParagraphIndex $d "I don't have experience with Power"
$frag = $pkgOpen + @'
<w:p>
  <w:r><w:t xml:space="preserve">I have researched about Power Apps and it seems an intuitive way of creating small apps for businesses. I </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>don't</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> have a project I've worked on in Power Apps, but I am willing to learn about the tool and other such technologies. </w:t></w:r>
</w:p>
'@ + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($frag)

# --- 2. "What is your weakness?" answer -------------------------------------
$idx = Find-ParagraphIndex $d "I would say one of my key weaknesses"
$frag = $pkgOpen + @'
<w:p>
  <w:r><w:t xml:space="preserve">One of my key weaknesses is taking tasks very seriously and passionately - which sometimes </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>isn't</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> reflected in other team members and creates communication gaps. I am working on it by regularly communicating with my teammates and understanding their point of view.  </w:t></w:r>
</w:p>
'@ + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($frag)

# --- 3. "Why do you want to work at Data-Driven?" heading + new answer ------
$idx = Find-ParagraphIndex $d "Why do you want to work at Data-Driven?"
$frag = $pkgOpen + @'
<w:p>
  <w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Why do you want to work at Data-Driven?</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">I have heard great comments about Data-Driven management and the work style from my seniors. This has really motivated me to work at the organization. Besides, I am also a </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>Masters of Data Science</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> &amp; AI student at UNSW and I believe Data-Driven would be able to provide me with the experience of working closely in projects related to my field. I am also interested in working with cloud technologies (and have some experience in it as well), which is something Data-Driven also follows so I believe I will be able to thrive in this environment.  </w:t></w:r>
</w:p>
'@ + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($frag)

# --- 4. Drop the old page-break hint on the visa-answer paragraph ----------
$idx = Find-ParagraphIndex $d "Currently, I am on a student visa"
$frag = $pkgOpen + @'
<w:p>
  <w:r><w:t xml:space="preserve">Currently, I am on a student visa (subclass 500) and I have the right to work 20 hours when the session is in place and unlimited during the holidays.  </w:t></w:r>
</w:p>
'@ + $pkgClose
$d.Paragraphs.Item($idx).Range.InsertXML($frag)

Write-Output "Applied Power BI Technical Test edits."
